$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column E (5): "Periodo Mora" labels, now sorted ascending (1801 .. 2003) for rows 16-42
$periods = @("1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812", `
             "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912", `
             "2001","2002","2003")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 5).Value = $periods[$i]
}

# Column F (6): "Salario Basico" - rows 16-23 -> 27578, rows 24-42 -> 31249
for ($r = 16; $r -le 23; $r++) {
    $ws.Cells.Item($r, 6).Value = 27578
}
for ($r = 24; $r -le 42; $r++) {
    $ws.Cells.Item($r, 6).Value = 31249
}

# Column G (7): "Valor Mora" - rows 16-42 -> 781242
for ($r = 16; $r -le 42; $r++) {
    $ws.Cells.Item($r, 7).Value = 781242
}
